$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "50.011"
$ws.Range("E1").Value = "30.011"

$ws.Range("D2").Value = "50.012"
$ws.Range("E2").Value = "30.012"

$ws.Range("D3").Value = "50.013"
$ws.Range("E3").Value = "30.013"

$ws.Range("P1").Value = "test11@test.com"
$ws.Range("P2").Value = "test21@test.com"
$ws.Range("P3").Value = "test31@test.com"
